# Add two new rows (119 and 120) of data to each of the 7 worksheets.
# Column A holds a date serial, column B holds a numeric amount.
# The new data, keyed by worksheet name (matches sheet1..sheet7 order):
#   row 119: date 45988, row 120: date 45989 (value 0)

$wb = $excel.ActiveWorkbook

$newData = @{
    "진양산업"   = 3202
    "넥스트아이" = 1093
    "삼보산업"   = 1313
    "YBM넷"      = 2000
    "NE능률"     = 712
    "위즈코프"   = 1524
    "대영포장"   = 2612
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($newData.ContainsKey($name)) {
        $val119 = $newData[$name]

        # Row 119
        $ws.Cells.Item(119, 1).Value = 45988
        $ws.Cells.Item(119, 2).Value = $val119

        # Row 120
        $ws.Cells.Item(120, 1).Value = 45989
        $ws.Cells.Item(120, 2).Value = 0

        # Match the date-style formatting used by the rest of column A (style index 2 == numFmtId 165)
        $ws.Cells.Item(119, 1).NumberFormat = $ws.Cells.Item(118, 1).NumberFormat
        $ws.Cells.Item(120, 1).NumberFormat = $ws.Cells.Item(118, 1).NumberFormat
    }
}
